# "argumento 'language' traducao de ataques"
#
# 1. Reword the instructional note in A1 (sheet1): the sentence describing
#    how the attack duration compares to the sum of per-prompt runtimes is
#    rephrased in Portuguese ("vai ser sempre superior" -> "total ... pode
#    ser superior").
# 2. Move the saved cell selection from A7 to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Duração do ataque, em segundos: `n(Atenção! A duração total do ataque pode ser superior à soma do 'runtime' de cada prompt)"

$ws.Range("A5").Select()
